$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet ("Sheet2") right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Put the research link in cell A1 of the new sheet
$ws2.Range("A1").Value = "http://gatherer.wizards.com/Pages/Search/Default.aspx?action=advanced&color=+![W]+![B]+![R]+![G]&text=+[creatures]+[you]+[control]"

# Move the selection on Sheet1 from the old G6 to D2
$ws1.Range("D2").Select()

# Make Sheet2 the active tab (matches activeTab="1" / tabSelected on sheet2)
$ws2.Activate()
